# Insert a new weekly data row at the top of the "Ají" price table (row 97),
# pushing the existing rows 97-151 down to 98-152, and populate the new
# row 97 with the latest week's price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 97:151 down by one row, carrying formatting with them.
$ws.Rows("97:97").Insert()

# Fill in the newly inserted row with the new weekly record.
$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44438
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = 100112021
$ws.Range("G97").Value = "Ají"
$ws.Range("H97").Value = "Inferno"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 25
$ws.Range("K97").Value = 32000
$ws.Range("L97").Value = 33000
$ws.Range("M97").Value = 32520
$ws.Range("N97").Value = "$/caja 12 kilos"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 2710
$ws.Range("Q97").Value = 12
$ws.Range("R97").Value = "Hortaliza"
